$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.010.62'
$ws.Range("E2").Value = '  -1.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.554.27'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.67'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3788'
$ws.Range("E7").Value = '  +2.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3237'
$ws.Range("E8").Value = '  -2.40%  '
$ws.Range("E9").Value = '  -13.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.123'
$ws.Range("E10").Value = '  -3.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07302'
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.32'
$ws.Range("E13").Value = '  -6.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.718'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.781'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.565.22'
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001086'
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06616'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.91'
$ws.Range("E19").Value = '  -3.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.412'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9984'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.95'
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.022.43'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.286'
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.529'
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.93'
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.856'
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.731.76'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.19'
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.106'
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.922'
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.652'
$ws.Range("E34").Value = '  -16.32%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08143'
$ws.Range("E35").Value = '  -1.97%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.226'
$ws.Range("E36").Value = '  -6.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.233'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06175'
$ws.Range("E38").Value = '  -3.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02278'
$ws.Range("E39").Value = '  -6.24%  '
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.218'
$ws.Range("E41").Value = '  -5.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.87'
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9994'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5931'
$ws.Range("E44").Value = '  -4.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.53'
$ws.Range("E45").Value = '  -3.59%  '
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5728'
$ws.Range("E47").Value = '  -4.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.933'
$ws.Range("E48").Value = '  -5.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.68'
$ws.Range("E49").Value = '  -3.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.157'
$ws.Range("E50").Value = '  -3.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06885'
$ws.Range("E51").Value = '  -4.23%  '
